$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.904399999999995
$ws.Range("C4").Value = -12.08619999999999
$ws.Range("B6").Value = 5.88
$ws.Range("B7").Value = 5.660399999999998
$ws.Range("B8").Value = 6.554799999999998
$ws.Range("C8").Value = -12.7546
$ws.Range("C9").Value = -10.3449
$ws.Range("C12").Value = -10.61269999999999
$ws.Range("B16").Value = 4.8447
$ws.Range("C17").Value = -14.59939999999999
$ws.Range("C18").Value = -12.6535
$ws.Range("C19").Value = -11.3091
$ws.Range("B20").Value = 9.575399999999989
$ws.Range("C20").Value = -12.5836
$ws.Range("B21").Value = 9.31829999999999
$ws.Range("C26").Value = -11.8334
$ws.Range("B28").Value = 5.988599999999999
$ws.Range("B29").Value = 5.533700000000005
$ws.Range("B30").Value = 5.461899999999998
$ws.Range("C31").Value = -13.0434
$ws.Range("B32").Value = 7.011499999999995
$ws.Range("C39").Value = -11.65630000000001
$ws.Range("B40").Value = 9.316899999999992
$ws.Range("C40").Value = -12.53740000000001
$ws.Range("C41").Value = -12.56050000000001
$ws.Range("C42").Value = -11.7169
$ws.Range("C43").Value = -12.87079999999999
$ws.Range("B46").Value = 5.685700000000002
$ws.Range("C47").Value = -12.01329999999999
$ws.Range("C48").Value = -12.99539999999999
$ws.Range("B51").Value = 5.551399999999997
$ws.Range("B52").Value = 5.610499999999996
$ws.Range("C54").Value = -12.4303
$ws.Range("B57").Value = 5.023399999999998
$ws.Range("B59").Value = 4.641800000000004
$ws.Range("B62").Value = 6.084699999999998
$ws.Range("C62").Value = -12.55600000000001
$ws.Range("C63").Value = -10.22259999999999
$ws.Range("C64").Value = -10.56769999999999
$ws.Range("B66").Value = 5.932799999999997
$ws.Range("B73").Value = 8.566999999999995
$ws.Range("B74").Value = 8.868999999999993
$ws.Range("C76").Value = -12.301
$ws.Range("B77").Value = 9.133000000000006
$ws.Range("C81").Value = -14.06249999999999
$ws.Range("C84").Value = -13.92309999999999
$ws.Range("C89").Value = -13.2407
$ws.Range("B92").Value = 4.706099999999998
$ws.Range("C94").Value = -10.8346
$ws.Range("B100").Value = 5.291500000000001
